$d = $word.ActiveDocument

$replacements = @(
    @("965×2=1930", "571×9=5139"),
    @("933×6=5598", "467×9=4203"),
    @("336×5=1680", "591×7=4137"),
    @("723×7=5061", "596×8=4768"),
    @("574×9=5166", "755×6=4530"),
    @("482×5=2410", "726×6=4356"),
    @("764×4=3056", "129×3=387"),
    @("337×6=2022", "956×2=1912"),
    @("707×6=4242", "175×5=875"),
    @("104×5=520",  "534×2=1068"),
    @("574×8=4592", "281×2=562"),
    @("788×2=1576", "606×8=4848"),
    @("360×6=2160", "570×2=1140"),
    @("800×9=7200", "327×9=2943"),
    @("181×7=1267", "165×4=660"),
    @("781×5=3905", "515×5=2575"),
    @("459×8=3672", "839×9=7551"),
    @("199×5=995",  "904×6=5424"),
    @("933×3=2799", "520×2=1040"),
    @("618×3=1854", "306×2=612"),
    @("121×6=726",  "234×3=702"),
    @("665×5=3325", "919×5=4595"),
    @("941×7=6587", "543×3=1629"),
    @("937×7=6559", "895×7=6265"),
    @("279×2=558",  "255×2=510")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
